# Add a new data row for the "mh-toc-parity" resource (maternal parity
# observation profile) under the existing "Observations" header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header row's formatting (border/fill/alignment) down into row 2
# first, so the new row matches the rest of the table (same style index as
# row 1) instead of picking up Excel's default cell style.
$ws.Range("A1:K1").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)

# Profile / Name / Code / Time Types / Value Types / Data Absent Reason.
# Category Code, Category VS, Code VS, Body Site and Method are left blank
# for this resource.
$ws.Range("A2").Value = "mh-toc-parity"
$ws.Range("B2").Value = "MH TOC Parity Profile"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "LOINC#11977-6"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "dateTimeĵ, Periodĵ, Timingĵ, instantĵ"
$ws.Range("H2").Value = "integerĵ, timeĵ, dateTimeĵ, Periodĵ"
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
